# Update the StructureDefinition metadata (Metadata worksheet) and the
# corresponding value inside the Elements worksheet grid, mirroring the
# IG's move from the ibm.com/Alvearie namespace to linuxforhealth.org,
# bumping the version, refreshing the publish date, and clearing the
# duplicated base-Extension constraint text.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/encrypted-state"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Extension.url's Fixed Value mirrors the StructureDefinition's own URL.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/encrypted-state"

# The base "Extension" row's Constraint(s) cell no longer repeats the
# ele-1/ext-1 constraint text (it now only shows up lower in the table).
$elements.Range("AI2").Value = ""
